$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1086503871299013
$ws.Range("C2").Value = 3.098190120925032
$ws.Range("D2").Value = 0.2102712380257867

$ws.Range("B3").Value = 0.06028007884650129
$ws.Range("C3").Value = 5.963886658390839
$ws.Range("D3").Value = 0.3213260644776185

$ws.Range("B4").Value = 0.08563621982942979
$ws.Range("C4").Value = 8.172386766084804
$ws.Range("D4").Value = 0.2674528070743238

$ws.Range("B5").Value = 0.1289400195320898
$ws.Range("C5").Value = 10.36018747393904
$ws.Range("D5").Value = 0.2841313876750753

$ws.Range("B6").Value = 0.08407321901701156
$ws.Range("C6").Value = 12.35255299513625
$ws.Range("D6").Value = 0.4158146589852765

$ws.Range("B7").Value = 0.08543735910689121
$ws.Range("C7").Value = 14.71842611129444
$ws.Range("D7").Value = 0.3283951127639914

$ws.Range("B8").Value = 0.09742130204431984
$ws.Range("C8").Value = 17.1971104515109
$ws.Range("D8").Value = 0.1788360477047702

$ws.Range("B9").Value = 0.1427313361246129
$ws.Range("C9").Value = 19.13140524515815
$ws.Range("D9").Value = 0.3324845719457836

$ws.Range("B10").Value = 0.1297823871938461
$ws.Range("C10").Value = 21.31821265264314
$ws.Range("D10").Value = 0.3089838900677563

$ws.Range("B11").Value = 0.1454933390831186
$ws.Range("C11").Value = 23.892298645281
$ws.Range("D11").Value = 0.1587934702353489

$ws.Range("B12").Value = 0.1430481615477018
$ws.Range("C12").Value = 25.66333229387858
$ws.Range("D12").Value = 0.4208864622457664

$ws.Range("B13").Value = 0.09442821049168963
$ws.Range("C13").Value = 28.49226396390035
$ws.Range("D13").Value = 0.4419550303532804

$ws.Range("B14").Value = 0.07010726278135324
$ws.Range("C14").Value = 30.72832861619367
$ws.Range("D14").Value = 0.3229970084316202

$ws.Range("B15").Value = 0.09355412441684832
$ws.Range("C15").Value = 33.0914762897524
$ws.Range("D15").Value = 0.3899676978770101

$ws.Range("B16").Value = 0.1405829833570224
$ws.Range("C16").Value = 35.05212170287883
$ws.Range("D16").Value = 0.2751651969038366

